$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value is ambiguous with a plain number (e.g. "7.99").
# Excel auto-converts Range.Value for such strings into a Double, which would
# silently strip formatting like trailing zeros ("8.00" -> 8) or change the
# stored type from text to number. Temporarily force a Text number format so
# the assignment keeps the exact string, then restore the original "Normal"
# cell style so no stray formatting is left behind.
$textForceCells = @(
    "D5",
    "D6",
    "D8",
    "D11",
    "D12",
    "D14",
    "D21",
    "D22",
    "D24",
    "D25",
    "D27",
    "D30",
    "D31",
    "D32",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49"
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated coin price (column D) and 1h volume change (column E) values.
$ws.Range('D2').Value = '60.673.92'
$ws.Range('E2').Value = '  -0.42%  '
$ws.Range('D3').Value = '2.364.91'
$ws.Range('E3').Value = '  -3.63%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '542.83'
$ws.Range('E5').Value = '  -1.07%  '
$ws.Range('D6').Value = '140.37'
$ws.Range('E6').Value = '  -3.46%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '0.540'
$ws.Range('E8').Value = '  -9.82%  '
$ws.Range('D9').Value = '2.366.03'
$ws.Range('E9').Value = '  -3.53%  '
$ws.Range('E10').Value = '  -2.37%  '
$ws.Range('D11').Value = '0.154'
$ws.Range('E11').Value = '  +0.50%  '
$ws.Range('D12').Value = '5.31'
$ws.Range('E12').Value = '  -1.57%  '
$ws.Range('E13').Value = '  -2.58%  '
$ws.Range('D14').Value = '25.38'
$ws.Range('E14').Value = '  -2.11%  '
$ws.Range('D15').Value = '2.792.59'
$ws.Range('E15').Value = '  -3.48%  '
$ws.Range('E16').Value = '  -1.18%  '
$ws.Range('D17').Value = '60.653.25'
$ws.Range('E17').Value = '  -0.33%  '
$ws.Range('D18').Value = '2.365.20'
$ws.Range('E18').Value = '  -3.52%  '
$ws.Range('E19').Value = '  -3.43%  '
$ws.Range('E20').Value = '  -1.51%  '
$ws.Range('D21').Value = '316.30'
$ws.Range('E21').Value = '  -0.75%  '
$ws.Range('D22').Value = '6.68'
$ws.Range('E22').Value = '  -3.25%  '
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('D24').Value = '1.85'
$ws.Range('E24').Value = '  +4.42%  '
$ws.Range('D25').Value = '62.91'
$ws.Range('E25').Value = '  -0.80%  '
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('D27').Value = '7.80'
$ws.Range('E27').Value = '  +3.05%  '
$ws.Range('D28').Value = '2.483.65'
$ws.Range('E28').Value = '  -3.42%  '
$ws.Range('D29').Value = '0.0₃0925'
$ws.Range('E29').Value = '  -4.96%  '
$ws.Range('D30').Value = '519.98'
$ws.Range('E30').Value = '  -3.71%  '
$ws.Range('D31').Value = '1.42'
$ws.Range('E31').Value = '  -4.48%  '
$ws.Range('D32').Value = '7.99'
$ws.Range('E32').Value = '  -4.11%  '
$ws.Range('E33').Value = '  -3.94%  '
$ws.Range('E34').Value = '  -3.28%  '
$ws.Range('E35').Value = '  +0.12%  '
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('D37').Value = '4.63'
$ws.Range('E37').Value = '  -4.48%  '
$ws.Range('D38').Value = '5.43'
$ws.Range('E38').Value = '  -6.66%  '
$ws.Range('D39').Value = '0.373'
$ws.Range('E39').Value = '  -0.68%  '
$ws.Range('D40').Value = '17.96'
$ws.Range('E40').Value = '  -2.71%  '
$ws.Range('D41').Value = '1.73'
$ws.Range('E41').Value = '  +0.70%  '
$ws.Range('E42').Value = '  +0.08%  '
$ws.Range('D43').Value = '137.05'
$ws.Range('E43').Value = '  -5.61%  '
$ws.Range('D44').Value = '40.16'
$ws.Range('E44').Value = '  +0.74%  '
$ws.Range('D45').Value = '2.20'
$ws.Range('E45').Value = '  -4.48%  '
$ws.Range('D46').Value = '139.23'
$ws.Range('E46').Value = '  -4.27%  '
$ws.Range('D47').Value = '3.54'
$ws.Range('E47').Value = '  -0.38%  '
$ws.Range('D48').Value = '20.23'
$ws.Range('E48').Value = '  -2.52%  '
$ws.Range('D49').Value = '0.0517'
$ws.Range('E49').Value = '  -2.40%  '
$ws.Range('E50').Value = '  -1.63%  '
$ws.Range('E51').Value = '  -3.42%  '

# Restore the default "Normal" style on the cells we forced to Text format,
# so their stored style matches the rest of the untouched text cells.
foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
